# Replace the "word" column (column B, rows 2-193) of Sheet1 with the new
# working set of words for this retrieval sequence.
#
# The new working-set list has 191 distinct words; the sheet still has 192
# data rows, so the final word is repeated once more for the last row (the
# shared-string table therefore still ends up with 191 unique word strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    "sperren",
    "wehen",
    "führen",
    "kümmern",
    "kehren",
    "gründen",
    "scheinen",
    "fischen",
    "treten",
    "irren",
    "siegen",
    "schnellen",
    "ächzen",
    "sprengen",
    "liefern",
    "schalten",
    "treffen",
    "drohen",
    "scheitern",
    "tauchen",
    "spielen",
    "streichen",
    "tropfen",
    "wärmen",
    "geben",
    "träumen",
    "greifen",
    "achten",
    "fügen",
    "ärgern",
    "brauchen",
    "leisten",
    "drehen",
    "loben",
    "morden",
    "wundern",
    "trennen",
    "ändern",
    "kranken",
    "messen",
    "folgen",
    "filmen",
    "binden",
    "spinnen",
    "formen",
    "sorgen",
    "deuten",
    "wachsen",
    "warnen",
    "tollen",
    "kosten",
    "wehtun",
    "landen",
    "segeln",
    "malen",
    "boxen",
    "klingen",
    "schlagen",
    "backen",
    "dringen",
    "bitten",
    "erben",
    "jubeln",
    "schenken",
    "stärken",
    "flüchten",
    "öffnen",
    "rasen",
    "werden",
    "weichen",
    "kichern",
    "saufen",
    "wetten",
    "gnaden",
    "orten",
    "betteln",
    "lockern",
    "zeigen",
    "plaudern",
    "räumen",
    "krachen",
    "knarren",
    "schlucken",
    "weigern",
    "zögern",
    "stecken",
    "freuen",
    "pflanzen",
    "äußern",
    "bauen",
    "passen",
    "streifen",
    "doppeln",
    "gelten",
    "stammen",
    "spüren",
    "bellen",
    "meinen",
    "bluten",
    "kriegen",
    "schwören",
    "sinken",
    "feuern",
    "platzen",
    "stehlen",
    "suchen",
    "heißen",
    "leeren",
    "wirken",
    "helfen",
    "mögen",
    "klettern",
    "proben",
    "feiern",
    "decken",
    "altern",
    "fließen",
    "heulen",
    "strahlen",
    "ziehen",
    "lesen",
    "wenden",
    "lügen",
    "tragen",
    "ehren",
    "posten",
    "leiden",
    "trauen",
    "zielen",
    "fangen",
    "bergen",
    "jagen",
    "schmecken",
    "fühlen",
    "werfen",
    "brauen",
    "eignen",
    "schwächen",
    "narren",
    "hauen",
    "schrecken",
    "heben",
    "warten",
    "lohnen",
    "pfeifen",
    "reiten",
    "planen",
    "töten",
    "füttern",
    "münzen",
    "hupen",
    "graben",
    "sterben",
    "grüßen",
    "mauern",
    "knien",
    "bremsen",
    "schulden",
    "hören",
    "seufzen",
    "schauen",
    "runden",
    "schreiten",
    "enden",
    "schützen",
    "fahren",
    "zeichnen",
    "melden",
    "quälen",
    "machen",
    "sichern",
    "biegen",
    "rufen",
    "fassen",
    "liegen",
    "schließen",
    "stimmen",
    "dienen",
    "reizen",
    "mühen",
    "heilen",
    "schwingen",
    "fallen",
    "locken",
    "treiben",
    "weinen",
    "arten",
    "flehen",
    "spannen",
    "zünden",
    "wüten"
)

$startRow = 2
$col = 2  # column B = "word"

for ($i = 0; $i -lt $newWords.Length; $i++) {
    $ws.Cells.Item($startRow + $i, $col).Value = $newWords[$i]
}

# The sheet had one more data row (193) than the new word list has entries;
# carry the final word down into that last row so every trial row keeps a
# word value.
$lastRow = $startRow + $newWords.Length
$ws.Cells.Item($lastRow, $col).Value = $newWords[$newWords.Length - 1]

